$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 78 (shifts existing rows 78..113 down to 79..114)
$ws.Rows.Item(78).Insert()

# Populate the newly inserted row 78 with the new record
$ws.Cells.Item(78, 1).Value = 6
$ws.Cells.Item(78, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(78, 3).Value = "Metropolitana"
$ws.Cells.Item(78, 4).Value = 44455
$ws.Cells.Item(78, 5).Value = 13
$ws.Cells.Item(78, 6).Value = 100112026
$ws.Cells.Item(78, 7).Value = "Haba"
$ws.Cells.Item(78, 8).Value = "Sin especificar"
$ws.Cells.Item(78, 9).Value = "Primera"
$ws.Cells.Item(78, 10).Value = 600
$ws.Cells.Item(78, 11).Value = 10000
$ws.Cells.Item(78, 12).Value = 12000
$ws.Cells.Item(78, 13).Value = 10767
$ws.Cells.Item(78, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(78, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(78, 16).Value = 431
$ws.Cells.Item(78, 17).Value = 25
$ws.Cells.Item(78, 18).Value = "Hortaliza"
